$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new tracking-number text. These look like numbers, so a
# plain Value assignment would coerce them to doubles; prefixing with an
# apostrophe forces Excel to store them as text (matching the workbook's
# existing convention of keeping these PackageTrackNum/ShipmentTrackNum
# columns as shared strings). Resetting the Style back to "Normal"
# afterwards strips the quote-prefix formatting so the cell keeps its
# original (default) style.
$updates = [ordered]@{
    "C2"  = "320018538422"
    "C3"  = "320018538433"
    "C4"  = "320018538466"
    "C5"  = "320018538488"
    "D5"  = "320018538488"
    "C6"  = "320018538525"
    "D6"  = "320018538525"
    "C7"  = "320018538547"
    "D7"  = "320018538547"
    "C8"  = "320018538570"
    "C9"  = "320018538591"
    "C10" = "320018538628"
    "C11" = "320018538640"
    "C12" = "320018538694"
    "C13" = "320018538710"
    "D13" = "320018538710"
    "C14" = "320018538742"
    "D14" = "320018538742"
    "C15" = "320018538775"
    "D15" = "320018538775"
    "C16" = "320018538801"
    "D16" = "320018538801"
    "C17" = "320018538823"
    "D17" = "320018538823"
    "C18" = "320018538867"
    "C19" = "320018538889"
    "C20" = "320018538915"
    "C21" = "320018538937"
    "C22" = "320018538960"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}

# Scenario Q3 moved from a PASS to a FAIL result.
$ws.Range("Q3").Value = "FAIL"
